# Updates cryptos list figures (price + 1h volume change) and
# re-orders a few coin rows, matching the upstream GitHub Actions data
# refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.228.50"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.843.69"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'240.70"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'0.6739"
$ws.Range("E6").Value = "  -1.49%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.07424"
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").Value = "'0.2952"
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("D10").Value = "'22.88"
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("D11").Value = "'0.07716"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "1.830.48"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "'5.007"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").Value = "'0.6707"
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("D15").Value = "'86.11"
$ws.Range("E15").Value = "  -1.49%  "
$ws.Range("D16").Value = "'6.130"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "29.187.61"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000008310"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").Value = "'228.53"
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").Value = "'12.51"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "'7.186"
$ws.Range("E22").Value = "  -2.96%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("D24").Value = "'160.65"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").Value = "'8.693"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("E26").Value = "  -3.42%  "
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("D28").Value = "'1.508"
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("E29").Value = "  -1.91%  "
$ws.Range("D30").Value = "'4.067"
$ws.Range("E30").Value = "  -1.98%  "
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").Value = "'0.05308"
$ws.Range("E32").Value = "  +2.91%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.7601"
$ws.Range("E33").Value = "  -0.64%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'1.875"
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("D35").Value = "'1.135"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").Value = "'2.676"
$ws.Range("D37").Value = "1.328.78"
$ws.Range("E37").Value = "  +1.84%  "
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("D40").Value = "'0.9185"
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("D41").Value = "'5.949"
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("D42").Value = "'1.002"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").Value = "'103.43"
$ws.Range("E43").Value = "  -1.03%  "
$ws.Range("D44").Value = "'0.08037"
$ws.Range("E44").Value = "  +15.64%  "
$ws.Range("D45").Value = "1.977.87"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000123"
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.5162"
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.774"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").Value = "'63.84"
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("D50").Value = "'9.144"
$ws.Range("E50").Value = "  -3.94%  "
$ws.Range("D51").Value = "'0.05950"
$ws.Range("E51").Value = "  +0.34%  "
